# Updates cryptocurrency price/volume figures (and re-ranks two coin pairs)
# on the active worksheet, matching the latest GitHub Actions data refresh.
# Price-column (D) values are prefixed with a leading apostrophe -- the
# same trick Excel's UI uses -- so numeric-looking strings like "212.46"
# or "0.502" stay stored as text instead of being auto-coerced to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''26.254.10'
$ws.Cells.Item(2, 5).Value = '  -0.19%  '
$ws.Cells.Item(3, 4).Value = '''1.592.05'
$ws.Cells.Item(3, 5).Value = '  +0.02%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '''212.46'
$ws.Cells.Item(5, 5).Value = '  +0.49%  '
$ws.Cells.Item(6, 4).Value = '''0.502'
$ws.Cells.Item(6, 5).Value = '  -0.68%  '
$ws.Cells.Item(7, 5).Value = '  -0.01%  '
$ws.Cells.Item(8, 4).Value = '''0.245'
$ws.Cells.Item(8, 5).Value = '  -0.59%  '
$ws.Cells.Item(9, 5).Value = '  -0.60%  '
$ws.Cells.Item(10, 4).Value = '''19.06'
$ws.Cells.Item(10, 5).Value = '  -1.69%  '
$ws.Cells.Item(11, 4).Value = '''0.0852'
$ws.Cells.Item(11, 5).Value = '  +0.50%  '
$ws.Cells.Item(12, 4).Value = '''1.815.61'
$ws.Cells.Item(12, 5).Value = '  +0.06%  '
$ws.Cells.Item(13, 4).Value = '''1.600.02'
$ws.Cells.Item(13, 5).Value = '  +0.85%  '
$ws.Cells.Item(14, 5).Value = '  -1.99%  '
$ws.Cells.Item(15, 5).Value = '  -2.25%  '
$ws.Cells.Item(16, 4).Value = '''63.82'
$ws.Cells.Item(16, 5).Value = '  -1.20%  '
$ws.Cells.Item(17, 4).Value = '''26.232.45'
$ws.Cells.Item(17, 5).Value = '  -0.27%  '
$ws.Cells.Item(18, 5).Value = '  -0.73%  '
$ws.Cells.Item(19, 4).Value = '''215.65'
$ws.Cells.Item(19, 5).Value = '  +1.59%  '
$ws.Cells.Item(20, 4).Value = '''7.28'
$ws.Cells.Item(20, 5).Value = '  -3.14%  '
$ws.Cells.Item(22, 4).Value = '''4.30'
$ws.Cells.Item(23, 4).Value = '''9.06'
$ws.Cells.Item(23, 5).Value = '  +0.58%  '
$ws.Cells.Item(24, 4).Value = '''2.13'
$ws.Cells.Item(24, 5).Value = '  -0.98%  '
$ws.Cells.Item(25, 4).Value = '''144.60'
$ws.Cells.Item(25, 5).Value = '  +0.62%  '
$ws.Cells.Item(26, 5).Value = '  -0.07%  '
$ws.Cells.Item(27, 4).Value = '''6.96'
$ws.Cells.Item(27, 5).Value = '  -1.52%  '
$ws.Cells.Item(28, 5).Value = '  -0.71%  '
$ws.Cells.Item(29, 4).Value = '''15.13'
$ws.Cells.Item(29, 5).Value = '  -0.72%  '
$ws.Cells.Item(30, 4).Value = '''0.0491'
$ws.Cells.Item(30, 5).Value = '  -2.43%  '
$ws.Cells.Item(31, 5).Value = '  +0.14%  '
$ws.Cells.Item(32, 5).Value = '  -0.87%  '
$ws.Cells.Item(33, 4).Value = '''1.422.99'
$ws.Cells.Item(33, 5).Value = '  +7.14%  '
$ws.Cells.Item(34, 5).Value = '  -1.38%  '
$ws.Cells.Item(35, 5).Value = '  -0.39%  '
$ws.Cells.Item(36, 4).Value = '''1.46'
$ws.Cells.Item(36, 5).Value = '  -0.83%  '
$ws.Cells.Item(37, 5).Value = '  -3.31%  '
$ws.Cells.Item(38, 4).Value = '''0.0165'
$ws.Cells.Item(38, 5).Value = '  -1.07%  '
$ws.Cells.Item(39, 2).Value = 'FraxShare'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(39, 4).Value = '''5.92'
$ws.Cells.Item(39, 5).Value = '  +3.75%  '
$ws.Cells.Item(40, 2).Value = 'ARBITRUM'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(40, 4).Value = '''0.823'
$ws.Cells.Item(40, 5).Value = '  +0.61%  '
$ws.Cells.Item(41, 5).Value = '  -0.04%  '
$ws.Cells.Item(42, 4).Value = '''0.989'
$ws.Cells.Item(42, 5).Value = '  -2.46%  '
$ws.Cells.Item(43, 2).Value = 'MXToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(43, 4).Value = '''2.14'
$ws.Cells.Item(43, 5).Value = '  +0.08%  '
$ws.Cells.Item(44, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(44, 4).Value = '''0.765'
$ws.Cells.Item(44, 5).Value = '  -0.08%  '
$ws.Cells.Item(45, 4).Value = '''1.727.93'
$ws.Cells.Item(45, 5).Value = '  +0.08%  '
$ws.Cells.Item(46, 4).Value = '''60.98'
$ws.Cells.Item(46, 5).Value = '  -1.59%  '
$ws.Cells.Item(47, 4).Value = '''86.68'
$ws.Cells.Item(47, 5).Value = '  -1.48%  '
$ws.Cells.Item(48, 4).Value = '''1.48'
$ws.Cells.Item(48, 5).Value = '  +0.06%  '
$ws.Cells.Item(49, 4).Value = '''0.0501'
$ws.Cells.Item(49, 5).Value = '  -0.66%  '
$ws.Cells.Item(50, 4).Value = '''0.0952'
$ws.Cells.Item(50, 5).Value = '  -2.69%  '
$ws.Cells.Item(51, 5).Value = '  -0.07%  '
